# Auto-generated edit script: adds 2023-08-29 data increments to violent-crime-full-year workbook
# Updates column J (year 2023 cumulative totals) and a few neighboring columns across 46 worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 4995
$ws.Range('I3').Value = 7489
$ws.Range('J3').Value = 5310
$ws.Range('E4').Value = 2002
$ws.Range('F4').Value = 1895
$ws.Range('G4').Value = 1469
$ws.Range('J4').Value = 1185
$ws.Range('J5').Value = 419
$ws.Range('I6').Value = 8964
$ws.Range('J6').Value = 6590
$ws.Range('E7').Value = 26007
$ws.Range('F7').Value = 24086
$ws.Range('G7').Value = 24693
$ws.Range('J7').Value = 18499

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J6').Value = 154
$ws.Range('J7').Value = 243

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 297
$ws.Range('J7').Value = 803

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J3').Value = 146
$ws.Range('J7').Value = 392

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J2').Value = 84
$ws.Range('J3').Value = 106
$ws.Range('J7').Value = 284

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 273
$ws.Range('J6').Value = 191
$ws.Range('J7').Value = 715

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J6').Value = 45
$ws.Range('J7').Value = 169

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J5').Value = 4
$ws.Range('J6').Value = 47
$ws.Range('J7').Value = 153

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J5').Value = 58
$ws.Range('J6').Value = 139
$ws.Range('J7').Value = 542
$ws.Range('J8').Value = 1176
$ws.Range('J10').Value = 125
$ws.Range('J11').Value = 284
$ws.Range('J14').Value = 86
$ws.Range('J19').Value = 533
$ws.Range('J20').Value = 386
$ws.Range('J21').Value = 48
$ws.Range('J24').Value = 54
$ws.Range('J25').Value = 91
$ws.Range('J29').Value = 1033
$ws.Range('J31').Value = 169
$ws.Range('J33').Value = 845
$ws.Range('J36').Value = 257
$ws.Range('J37').Value = 580
$ws.Range('J41').Value = 121
$ws.Range('J42').Value = 749
$ws.Range('J44').Value = 139
$ws.Range('J46').Value = 64
$ws.Range('I47').Value = 188
$ws.Range('J48').Value = 210
$ws.Range('J50').Value = 110
$ws.Range('J51').Value = 236
$ws.Range('J52').Value = 471
$ws.Range('J53').Value = 243
$ws.Range('J54').Value = 354
$ws.Range('J55').Value = 233
$ws.Range('E63').Value = 346
$ws.Range('F63').Value = 185
$ws.Range('G63').Value = 268
$ws.Range('I63').Value = 236
$ws.Range('J63').Value = 68
$ws.Range('J65').Value = 485
$ws.Range('J67').Value = 715
$ws.Range('J76').Value = 269
$ws.Range('J77').Value = 149
$ws.Range('J78').Value = 230
$ws.Range('J79').Value = 533
$ws.Range('J83').Value = 392
$ws.Range('J84').Value = 153
$ws.Range('J85').Value = 803
$ws.Range('J86').Value = 116
$ws.Range('J88').Value = 203
$ws.Range('J91').Value = 205
$ws.Range('J97').Value = 148
$ws.Range('J98').Value = 122
$ws.Range('J99').Value = 284
$ws.Range('E101').Value = 26007
$ws.Range('F101').Value = 24086
$ws.Range('G101').Value = 24693
$ws.Range('J101').Value = 18499

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 202
$ws.Range('J7').Value = 580

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 211
$ws.Range('J3').Value = 276
$ws.Range('J5').Value = 37
$ws.Range('J7').Value = 845

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 143
$ws.Range('J4').Value = 20
$ws.Range('J7').Value = 485

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J6').Value = 165
$ws.Range('J7').Value = 354

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J4').Value = 58
$ws.Range('J6').Value = 269
$ws.Range('J7').Value = 1033

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 131
$ws.Range('J3').Value = 154
$ws.Range('J6').Value = 198
$ws.Range('J7').Value = 533

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J6').Value = 50
$ws.Range('J7').Value = 139

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J4').Value = 34
$ws.Range('J7').Value = 210

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J6').Value = 146
$ws.Range('J7').Value = 269

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 171
$ws.Range('J5').Value = 14
$ws.Range('J6').Value = 172
$ws.Range('J7').Value = 542

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 51
$ws.Range('J7').Value = 139

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J2').Value = 26
$ws.Range('J7').Value = 121

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 164
$ws.Range('J3').Value = 150
$ws.Range('J6').Value = 383
$ws.Range('J7').Value = 749

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J2').Value = 28
$ws.Range('J7').Value = 125

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J2').Value = 63
$ws.Range('J3').Value = 76
$ws.Range('J6').Value = 62
$ws.Range('J7').Value = 230

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J2').Value = 55
$ws.Range('J7').Value = 233

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J6').Value = 12
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('J6').Value = 26
$ws.Range('J7').Value = 64

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J2').Value = 63
$ws.Range('J7').Value = 205

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('J6').Value = 30
$ws.Range('J7').Value = 48

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J5').Value = 14
$ws.Range('J7').Value = 533

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 105
$ws.Range('J7').Value = 386

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J3').Value = 80
$ws.Range('J6').Value = 73
$ws.Range('J7').Value = 257

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J3').Value = 58
$ws.Range('J7').Value = 284

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J2').Value = 40
$ws.Range('J7').Value = 91

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('I6').Value = 60
$ws.Range('I7').Value = 188

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J6').Value = 74
$ws.Range('J7').Value = 122

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J6').Value = 33
$ws.Range('J7').Value = 110

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 110
$ws.Range('J3').Value = 139
$ws.Range('J7').Value = 471

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J6').Value = 100
$ws.Range('J7').Value = 148

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J6').Value = 90
$ws.Range('J7').Value = 203

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J6').Value = 27
$ws.Range('J7').Value = 86

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J2').Value = 20
$ws.Range('J7').Value = 58

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J2').Value = 15
$ws.Range('J4').Value = 61
$ws.Range('J6').Value = 22
$ws.Range('J7').Value = 116

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 65
$ws.Range('J7').Value = 236

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 53
$ws.Range('J3').Value = 53
$ws.Range('J7').Value = 149

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J4').Value = 71
$ws.Range('J6').Value = 388
$ws.Range('J7').Value = 1176
